# Applies the "edits and sim dir" commit:
#  - view changes: zoom to 210%, active selection moves to D12
#  - parameter-sheet value tweaks across several rows, including
#    swapping two static D4/D5 formulas for literal values and
#    turning the D10/D11 literals into geometric-mean formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (k) ---
$ws.Range("E2").Value = 0.0001

# --- Row 3 (sy) ---
$ws.Range("E3").Value = 0.05

# --- Row 4 (cghbw) ---
$ws.Range("D4").Value = 0.1
$ws.Range("E4").Value = 0.001
$ws.Range("F4").Value = 1

# --- Row 5 (cghbe) ---
$ws.Range("D5").Value = 0.1
$ws.Range("E5").Value = 0.001
$ws.Range("F5").Value = 1

# --- Row 8 (criv) ---
$ws.Range("E8").Value = 0.000001
$ws.Range("F8").Value = 0.01

# --- Row 9 (cdrn) ---
$ws.Range("E9").Value = 0.0001
$ws.Range("F9").Value = 0.01

# --- Row 10 (tsat): update bounds first, then derive D10 via formula ---
$ws.Range("F10").Value = 0.5
$ws.Range("D10").Formula = "=10^((LOG10(E10)+LOG10(F10))/2)"

# --- Row 11 (dmax): update bounds first, then derive D11 via formula ---
$ws.Range("F11").Value = 50
$ws.Range("D11").Formula = "=10^((LOG10(E11)+LOG10(F11))/2)"

# Recalculate so every formula cell carries a fresh cached value.
$excel.CalculateFull()

# --- View state: zoom + active cell/selection ---
$excel.ActiveWindow.Zoom = 210
$ws.Range("D12").Select()
